$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("elisabeth venchiarutti", 4, 3, $null, $null, $null, $null, 1, $null, $null),
    @("hello again", 3, 1, 1, 1, 1, 1, $null, $null, $null),
    @("Victoire", 3, 2, $null, 1, 1, 1, $null, $null, $null),
    @("Victoire", 5, 3, $null, $null, $null, $null, $null, $null, $null),
    @("Le grand succès", 1, 2, 4, 1, $null, $null, $null, $null, $null),
    @("Marc Clément ", 5, $null, $null, 1, $null, $null, $null, 2, $null),
    @("Elisabeth jeudi 15 aout", 2, 1, 3, $null, $null, 1, 1, $null, $null),
    @("hourr ahourra", 3, 1, 1, $null, $null, 1, 1, 1, $null),
    @("Nathalie Marcot", 1, 2, $null, 2, $null, $null, $null, 1, 2),
    @("Nouvelle version", 2, 6, $null, $null, $null, $null, $null, $null, $null),
    @("cest moi", 5, 3, $null, $null, $null, $null, $null, $null, $null)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($row, $c + 1).Value = $val
        }
    }
}
